# Add season-record columns (Wins/Losses/Ties) to the roster/stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers, styled like the rest of the header row ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font + border) from an existing header cell
# onto the three new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-47): the team finished the season 87-75-0, so stamp the
# same season record on every player's row. ---
$lastRow = 47
$firstDataRow = 2
$rowCount = $lastRow - $firstDataRow + 1

$data = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = 87
    $data[$i,1] = 75
    $data[$i,2] = 0
}

$ws.Range("AD" + $firstDataRow + ":AF" + $lastRow).Value = $data

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
